$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.000.58"
$ws.Range("E2").Value = "  -0.39%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.236.34"
$ws.Range("E3").Value = "  -0.66%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.65"
$ws.Range("E5").Value = "  -4.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.40"
$ws.Range("E6").Value = "  -6.87%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.570"
$ws.Range("E7").Value = "  -0.74%  "

$ws.Range("E8").Value = "  +0.26%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.522"
$ws.Range("E9").Value = "  -4.32%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.75"
$ws.Range("E10").Value = "  -6.01%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0808"
$ws.Range("E11").Value = "  -2.44%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.19"
$ws.Range("E12").Value = "  -4.62%  "

$ws.Range("E13").Value = "  -0.95%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.577.36"
$ws.Range("E14").Value = "  -0.59%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.234.29"
$ws.Range("E15").Value = "  -2.10%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.821"
$ws.Range("E16").Value = "  -3.43%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.52"
$ws.Range("E17").Value = "  -5.02%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.834.46"
$ws.Range("E18").Value = "  -0.59%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0962"
$ws.Range("E19").Value = "  -1.66%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.05"
$ws.Range("E20").Value = "  -9.89%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.28"
$ws.Range("E21").Value = "  -2.53%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.38"
$ws.Range("E22").Value = "  -0.33%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.85"
$ws.Range("E23").Value = "  +0.62%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.91"
$ws.Range("E24").Value = "  -5.86%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.96"
$ws.Range("E25").Value = "  -4.92%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.01"
$ws.Range("E26").Value = "  +0.28%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.87"
$ws.Range("E27").Value = "  -6.08%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.19"
$ws.Range("E28").Value = "  -1.19%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "37.68"
$ws.Range("E29").Value = "  -0.39%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.99"
$ws.Range("E30").Value = "  -2.04%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.91"
$ws.Range("E31").Value = "  -1.01%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "150.92"
$ws.Range("E32").Value = "  -5.00%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0797"
$ws.Range("E33").Value = "  -5.94%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").Value = "3.20"
$ws.Range("E34").Value = "  +0.22%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").Value = "2.58"
$ws.Range("E35").Value = "  -3.42%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.110"
$ws.Range("E36").Value = "  -2.07%  "

$ws.Range("E37").Value = "  +0.86%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.78"
$ws.Range("E38").Value = "  -8.79%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.04"
$ws.Range("E39").Value = "  -6.77%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.84"
$ws.Range("E40").Value = "  -7.78%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.35"
$ws.Range("E41").Value = "  -9.05%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0295"
$ws.Range("E42").Value = "  -6.25%  "

$ws.Range("E43").Value = "  +0.20%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.730.67"
$ws.Range("E44").Value = "  -0.79%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "84.68"
$ws.Range("E45").Value = "  +2.56%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.187"
$ws.Range("E46").Value = "  -5.33%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "99.75"
$ws.Range("E47").Value = "  -2.95%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.92"
$ws.Range("E48").Value = "  -4.72%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.05"
$ws.Range("E49").Value = "  -2.72%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("B50").Value = "ordi"
$ws.Range("C50").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D50").Value = "68.58"
$ws.Range("E50").Value = "  -8.23%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").Value = "53.86"
$ws.Range("E51").Value = "  -6.66%  "
